$d = $word.ActiveDocument

# --- 1. "a day-to-day life" -> "day-to-day life" (drop redundant "a ") ---
$d.Content.Find.Execute("a day-to-day life", $true, $false, $false, $false, $false, $true, 1, $false, "day-to-day life", 2) | Out-Null

# --- 2. "in terms of usability they lack" -> "in terms of usability, they lack" ---
$d.Content.Find.Execute("in terms of usability they lack", $true, $false, $false, $false, $false, $true, 1, $false, "in terms of usability, they lack", 2) | Out-Null

# --- 3. "...obvious reasons like to ensure the bug priority, increase coverage area, etc. On the other hand" ---
#        -> "...obvious reasons such as to ensure the bug priority, increase coverage area. On the other hand"
$d.Content.Find.Execute("obvious reasons like to ensure the bug priority, increase coverage area, etc. On the other hand", $true, $false, $false, $false, $false, $true, 1, $false, "obvious reasons such as to ensure the bug priority, increase coverage area. On the other hand", 2) | Out-Null

# --- 4. "This shows the future scope" -> "These reasons show the future scope" ---
$d.Content.Find.Execute("This shows the future scope", $true, $false, $false, $false, $false, $true, 1, $false, "These reasons show the future scope", 2) | Out-Null

# --- 5. Remove the stray empty run (with a bogus xhtml xmlns on its <w:t>) that sits between
#        "...is ongoing?" and the line break introducing "3. How to carry...". It carries no
#        text of its own, so locate it via the line-break character immediately following
#        "ongoing?" and collapse it away with a scoped, in-place Find/Replace on that one
#        character (scoping to a narrow Range keeps the rest of the document's breaks intact).
$probe = $d.Content
$probe.Find.Execute("ongoing?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$breakPos = $probe.End
$vt = [char]11
$breakRange = $d.Range($breakPos, $breakPos + 1)
$breakRange.Find.Execute($vt, $false, $false, $false, $false, $false, $true, 1, $false, $vt, 2) | Out-Null

# --- 6. Move the "_GoBack" bookmark from the end of the last paragraph into the middle of
#        "works" ("2. What feedback wo|rks to know..."). Adding a same-named bookmark
#        automatically replaces the previous one, wherever it was.
$target = $d.Content
$target.Find.Execute("2. What feedback wo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $target.End
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null
